$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restyle row 30 to the "closing" border style (s=6/7), matching row 3/6/... groups ---
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A30:E30").PasteSpecial(-4122) | Out-Null

# --- Step 2: seed rows 31-32 with the "open" style (s=4/5) by copying from an existing pair (row29 header, row5 continuation) ---
$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A31:E31").PasteSpecial(-4122) | Out-Null
$ws.Range("B5:E5").Copy() | Out-Null
$ws.Range("B32:E32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 3: row heights ---
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 21.6

# --- Step 4: numeric (non-shared-string) cells first, order does not affect sharedStrings table ---
$ws.Cells.Item(31, 2).Value = 19
$ws.Cells.Item(32, 2).Value = 22

# --- Step 5: text cells, written in the exact order the new strings must land in sharedStrings.xml
#            (108..114): C31, C32, A31, D32, D31, E31, E32 ---
$ws.Cells.Item(31, 3).Value = ' My partner [CS:N]Murkrow[CR] won\''t do\nanything I say.'   # C31 -> si 108
$ws.Cells.Item(32, 3).Value = ' I really want to go look for the\n[CS:I]Secret Slab[CR], but…'   # C32 -> si 109
$ws.Cells.Item(31, 1).Value = 'SCRIPT/T01P01A/us2304.ssb'   # A31 -> si 110
$ws.Cells.Item(32, 4).Value = ' Я очень хочу отправиться на\nпоиски [CS:I]Таблички-Секрета[CR], но...'   # D32 -> si 111
$ws.Cells.Item(31, 4).Value = ' Моя спутница [CS:N]Маркроу[CR] меня\nсовсем не слушает.'   # D31 -> si 112
$ws.Cells.Item(31, 5).Value = ' Íïÿ òðôóîéøà [CS:N]Íàñëñïô[CR] íåîÿ\nòïâòåí îå òìôšàåó.'   # E31 -> si 113
$ws.Cells.Item(32, 5).Value = ' Ÿ ïœåîû öïœô ïóðñàâéóûòÿ îà\nðïéòëé [CS:I]Óàáìéœëé-Òåëñåóà[CR], îï...'   # E32 -> si 114

# --- Step 6: selection, to mirror the diff's <selection activeCell="C29" sqref="C29"/> ---
$ws.Range("C29").Select() | Out-Null

Write-Output "done"
